$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.658.90"
$ws.Range("E2").Value = "  +4.43%  "
$ws.Range("D3").Value = "1.606.52"
$ws.Range("E3").Value = "  +3.75%  "
$ws.Range("E4").Value = "  -0.45%  "
$ws.Range("D5").Value = "'214.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.98%  "
$ws.Range("E6").Value = "  +8.04%  "
$ws.Range("D7").Value = "'0.996"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.46%  "
$ws.Range("D8").Value = "'26.78"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +11.53%  "
$ws.Range("E9").Value = "  +3.57%  "
$ws.Range("E10").Value = "  +3.06%  "
$ws.Range("D11").Value = "'0.0917"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.29%  "
$ws.Range("D12").Value = "1.835.03"
$ws.Range("E12").Value = "  +3.56%  "
$ws.Range("D13").Value = "1.616.77"
$ws.Range("E13").Value = "  +4.27%  "
$ws.Range("D14").Value = "29.677.57"
$ws.Range("E14").Value = "  +4.60%  "
$ws.Range("E15").Value = "  +3.98%  "
$ws.Range("D16").Value = "'0.529"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.09%  "
$ws.Range("D17").Value = "'246.62"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +7.72%  "
$ws.Range("D18").Value = "'63.74"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.78%  "
$ws.Range("E19").Value = "  +3.87%  "
$ws.Range("E20").Value = "  +3.70%  "
$ws.Range("D21").Value = "'0.996"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.45%  "
$ws.Range("E22").Value = "  +4.61%  "
$ws.Range("E23").Value = "  +4.25%  "
$ws.Range("E24").Value = "  +4.96%  "
$ws.Range("D25").Value = "'155.69"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.00%  "
$ws.Range("D26").Value = "'15.39"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.33%  "
$ws.Range("E27").Value = "  +6.41%  "
$ws.Range("E28").Value = "  +2.99%  "
$ws.Range("E29").Value = "  -0.46%  "
$ws.Range("D30").Value = "'0.0473"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.39%  "
$ws.Range("E31").Value = "  +0.65%  "
$ws.Range("E32").Value = "  +2.91%  "
$ws.Range("D33").Value = "1.442.37"
$ws.Range("E33").Value = "  +4.44%  "
$ws.Range("D34").Value = "'3.11"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.02%  "
$ws.Range("D35").Value = "'1.06"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.80%  "
$ws.Range("E36").Value = "  +11.26%  "
$ws.Range("E37").Value = "  +3.41%  "
$ws.Range("D38").Value = "'2.30"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.27%  "
$ws.Range("E39").Value = "  +3.01%  "
$ws.Range("B40").Value = "BitcoinSV"
$ws.Range("C40").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D40").Value = "'56.92"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +29.51%  "
$ws.Range("B41").Value = "ImmutableX"
$ws.Range("C41").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D41").Value = "'0.537"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.75%  "
$ws.Range("D42").Value = "'1.96"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.26%  "
$ws.Range("E43").Value = "  +4.42%  "
$ws.Range("D44").Value = "'0.996"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.42%  "
$ws.Range("D45").Value = "'67.37"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +8.80%  "
$ws.Range("D46").Value = "'0.0467"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.72%  "
$ws.Range("E47").Value = "  -0.12%  "
$ws.Range("D48").Value = "1.746.13"
$ws.Range("E48").Value = "  +3.68%  "
$ws.Range("D49").Value = "'86.22"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.77%  "
$ws.Range("D50").Value = "'0.838"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.96%  "
$ws.Range("D51").Value = "0.0₆0104"
$ws.Range("E51").Value = "  +1.34%  "
